# Update cryptocurrency price/volume data (refresh scrape results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.799.59"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.630.93"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.98"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0631"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.857.31"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "1.636.98"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.551"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "0.0₃0759"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.70"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "25.803.49"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.93"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.46"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "1.140.26"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.53"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.798"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "1.768.62"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "0.0₆0110"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.41"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("E49").Value = "  +6.07%  "
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.415"
$ws.Range("E51").Value = "  -0.70%  "


Write-Host "Applied cryptos list refresh."
